$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 111.0
$ws.Range("C1").Value = 1110.0
$ws.Range("B3").Value = 31.0
$ws.Range("C3").Value = 465.0
